# Portugal Segunda Liga - swap duplicated/misordered match rows.
# The source feed occasionally emitted a pair of adjacent rows with the
# columns B..AD swapped between them (column A, the running index, stays
# put). This script corrects six such pairs by swapping every column
# from B through AD between the two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) through AD (30)
$firstCol = 2
$lastCol  = 30

$rowPairs = @(
    @(39, 40),
    @(88, 89),
    @(140, 141),
    @(186, 187),
    @(243, 244),
    @(260, 261)
)

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
